# Version 2.0.1 solucionado error espera de base de datos
# Fill in the patient intake form with real data, replacing the placeholder
# values (22222222 / null / blank spaces) that were in the template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Patient name (Apellidos / Nombres) - row 6
$ws.Range("A6").Value = "CASTRO"
$ws.Range("C6").Value = "SANCHEZ"
$ws.Range("E6").Value = "DAISY"
$ws.Range("G6").Value = "CATALINA"
$ws.Range("I6").Value = "17-430/201758293"

# Direccion actual - row 8
$ws.Range("A8").Value = "20 AV 21-18 Z, 6 SAN JUAN DE DIOS"
$ws.Range("D8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("H8").Value = ""
$ws.Range("J8").Value = ""

# Direccion habitual - row 10
$ws.Range("A10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("H10").Value = ""
$ws.Range("J10").Value = ""

# Fecha nacimiento / edad / lugar nacimiento / sexo - row 12
$ws.Range("A12").Value = "4-6-80"
$ws.Range("F12").Value = "37"
$ws.Range("H12").Value = "HONDUREÑA"
$ws.Range("J12").Value = "FEMENINO"

# Estado civil / ocupacion / nacionalidad / cedula - row 14
$ws.Range("A14").Value = "Casado"
$ws.Range("D14").Value = "AMA DE CASA"
$ws.Range("F14").Value = "HONDUREÑA"
$ws.Range("H14").Value = "1804-198001758"

# Nombre del conyugue - row 16
$ws.Range("A16").Value = "JAVIER PLATA"
$ws.Range("F16").Value = ""

# Nombre del padre / madre - row 18
$ws.Range("A18").Value = "LUIS CASTRO"
$ws.Range("F18").Value = "MARIA ELENA SANCHEZ"

# Parentesco / direccion / telefono de emergencia - row 20
$ws.Range("A20").Value = ""
$ws.Range("F20").Value = "ESPOSO"
$ws.Range("H20").Value = ""
$ws.Range("J20").Value = "56106261"

# Otras hospitalizaciones / referido de - row 22
$ws.Range("A22").Value = ""
$ws.Range("F22").Value = ""

# Fecha de ingreso / hora / servicio - row 24
$ws.Range("A24").Value = "18-10-17"
$ws.Range("C24").Value = "08-41"
$ws.Range("D24").Value = "CL. 35"
